$d = $word.ActiveDocument

# --- Edit 1: rewrite the "your mom's sister" sentence fragment -------------
$d.Content.Find.Execute(
    "your mom’s sister lives in the area as well",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "a long time ago I made an arrangement with your mom’s sister",
    2
)

# --- Edit 2: split "Sincerely yours, Dad." into two paragraphs -------------
# Locate the paragraph that still holds the full "Sincerely yours, Dad."
# line and remember its 1-based index inside $d.Paragraphs.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Sincerely yours, Dad.*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $p = $d.Paragraphs($targetIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    # "Lilith’s Dad: Sincerely yours, Dad." -> split right after the comma,
    # i.e. keep "Lilith’s Dad: Sincerely yours," in the first paragraph and
    # " Dad." moves into a newly created second paragraph.
    $splitOffset = 30

    # Turn the trailing " Dad." (before the paragraph mark) into its own
    # paragraph by inserting a paragraph break right before it.
    $tail = $d.Range($pStart + $splitOffset, $pEnd - 1)
    $tail.InsertParagraphAfter()

    # Remove " Dad." from the end of the first paragraph.
    $tail = $d.Range($pStart + $splitOffset, $pEnd - 1)
    $tail.Text = ""

    # Fill the freshly-created (empty) paragraph with its own full line.
    $newPara = $d.Paragraphs($targetIndex + 1)
    $newPara.Range.InsertBefore("Lilith’s Dad: Dad.")
}

Write-Output "done"
